$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder ".." text that was entered in column F (rows 2-8),
# leaving the cell styling (right-aligned fill) in place but with no value.
$ws.Range("F2:F8").ClearContents()

# Move the active selection to F8, matching where the user last clicked.
$ws.Range("F8").Select()
